# Sync attendance_reports: rotate the "Recorded By" (column G) list of
# contributors so that the first entry moves to the end, e.g.
#   "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
# This is applied to every row that has more than one comma-separated
# value in column G, except the specific value "admin@admin.com, System"
# which is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    if ($val -eq "admin@admin.com, System") {
        continue
    }

    if ($val -like "*,*") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
            $cell.Value2 = $rotated
        }
    }
}
